$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the 12 "TestcodewordN" label rows with 10 single-letter labels (A..J)
$letters = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $letters.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $letters[$i]
}

# The data set shrank from 12 rows to 10 rows - remove the trailing two rows
$ws.Range("A11:A12").EntireRow.Delete()

# Update the active selection from G8 to E8
$ws.Range("E8").Select()
